$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held 20 data rows (timestamps 0-1900, rows 2-21).
# It now holds 30 data rows (timestamps 0-2900, rows 2-31): two brand-new rows
# were inserted at the top (now rows 2-3), the original 20 rows were pushed down
# to rows 4-23, and eight brand-new rows were appended at the bottom (rows 24-31).

# Work from the bottom up so that source rows are not overwritten before they
# are copied down to their new location.
$cols = @("C", "D", "E", "F", "G", "H")
for ($oldRow = 21; $oldRow -ge 2; $oldRow--) {
    $newRow = $oldRow + 2
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $ws.Range("$col$oldRow").Value2
    }
}

# Re-number the timestamp (A) and label (B) columns for the full 2-31 row range;
# these simply continue the existing 0,100,200,... / "struggle" pattern.
for ($row = 2; $row -le 31; $row++) {
    $ws.Range("A$row").Value = ($row - 2) * 100
    $ws.Range("B$row").Value = "struggle"
}

# Populate the two newly inserted rows at the top with their sensor data.
$ws.Range("C2").Value = -0.1967945098876953
$ws.Range("D2").Value = -0.0881298780441284
$ws.Range("E2").Value = -0.4354097247123718
$ws.Range("F2").Value = 0.0100699262883591
$ws.Range("G2").Value = -0.003861541194575109
$ws.Range("H2").Value = -0.02743906991518268
$ws.Range("C3").Value = 0.3726930618286133
$ws.Range("D3").Value = 0.3928739428520202
$ws.Range("E3").Value = -0.1955753564834594
$ws.Range("F3").Value = -0.04216528505238931
$ws.Range("G3").Value = -0.05587235412427344
$ws.Range("H3").Value = -0.005946585338334982

# Populate the eight newly appended rows at the bottom with their sensor data.
$ws.Range("C24").Value = 0.6416101455688477
$ws.Range("D24").Value = 0.1077315807342529
$ws.Range("E24").Value = -3.572567462921143
$ws.Range("F24").Value = 1.138076220239913
$ws.Range("G24").Value = 4.945667840996566
$ws.Range("H24").Value = -0.5911586260309147
$ws.Range("C25").Value = 1.180892944335938
$ws.Range("D25").Value = -0.3624088764190674
$ws.Range("E25").Value = 1.944910764694214
$ws.Range("F25").Value = 0.5868015289306701
$ws.Range("G25").Value = 4.023616756711703
$ws.Range("H25").Value = 0.7568838426044971
$ws.Range("C26").Value = -0.6099348068237305
$ws.Range("D26").Value = -0.0995303392410278
$ws.Range("E26").Value = 1.559979677200317
$ws.Range("F26").Value = -0.3439888250538894
$ws.Range("G26").Value = 1.417871174155451
$ws.Range("H26").Value = 1.152574896812441
$ws.Range("C27").Value = 0.7382268905639648
$ws.Range("D27").Value = 0.5965696573257446
$ws.Range("E27").Value = 0.3601601719856262
$ws.Range("F27").Value = -0.006142936684953748
$ws.Range("G27").Value = 0.247877272416142
$ws.Range("H27").Value = -0.5417658090591317
$ws.Range("C28").Value = 0.1256790161132812
$ws.Range("D28").Value = 0.4359270334243774
$ws.Range("E28").Value = -0.5883067846298218
$ws.Range("F28").Value = -0.1216962014001841
$ws.Range("G28").Value = -0.5174121899264155
$ws.Range("H28").Value = -0.2046180449578262
$ws.Range("C29").Value = 0.1187114715576171
$ws.Range("D29").Value = 0.2241333723068237
$ws.Range("E29").Value = -0.9467962980270386
$ws.Range("F29").Value = -0.007408298704090516
$ws.Range("G29").Value = -0.9081197368855365
$ws.Range("H29").Value = -0.06973525623277696
$ws.Range("C30").Value = 0.1263256072998047
$ws.Range("D30").Value = 0.5689128637313843
$ws.Range("E30").Value = -0.7026804089546204
$ws.Range("F30").Value = 0.06267290592801897
$ws.Range("G30").Value = -0.8978846316434917
$ws.Range("H30").Value = -0.05514305708359769
$ws.Range("C31").Value = -0.1413173675537109
$ws.Range("D31").Value = 0.4839025735855102
$ws.Range("E31").Value = -0.0290583968162536
$ws.Range("F31").Value = -0.03725966301803718
$ws.Range("G31").Value = -0.5934867311497137
$ws.Range("H31").Value = 0.0119181060973484
